$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.685.42'
$ws.Range("E2").Value = '  +0.53%  '

$ws.Range("D3").Value = '1.639.36'
$ws.Range("E3").Value = '  -0.43%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.65'
$ws.Range("E5").Value = '  +0.00%  '

$ws.Range("E6").Value = '  -1.31%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.06'
$ws.Range("E8").Value = '  -2.13%  '

$ws.Range("E9").Value = '  +0.13%  '

$ws.Range("E10").Value = '  -0.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0894'
$ws.Range("E11").Value = '  +0.27%  '

$ws.Range("D12").Value = '1.870.99'
$ws.Range("E12").Value = '  -0.45%  '

$ws.Range("D13").Value = '1.632.94'
$ws.Range("E13").Value = '  -0.83%  '

$ws.Range("E14").Value = '  -0.05%  '

$ws.Range("E15").Value = '  -5.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.64'
$ws.Range("E16").Value = '  +0.13%  '

$ws.Range("D17").Value = '27.672.97'
$ws.Range("E17").Value = '  +0.62%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.74'
$ws.Range("E18").Value = '  -0.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.75'
$ws.Range("E19").Value = '  +2.71%  '

$ws.Range("D20").Value = '0.0₃0722'
$ws.Range("E20").Value = '  -0.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  +0.01%  '

$ws.Range("E22").Value = '  -0.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.26'
$ws.Range("E23").Value = '  +4.83%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.03'
$ws.Range("E24").Value = '  +0.42%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.90'
$ws.Range("E25").Value = '  +1.59%  '

$ws.Range("E26").Value = '  -1.21%  '

$ws.Range("E27").Value = '  -1.76%  '

$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.19%  '

$ws.Range("E30").Value = '  +0.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0488'
$ws.Range("E31").Value = '  +0.14%  '

$ws.Range("D33").Value = '1.457.46'
$ws.Range("E33").Value = '  +2.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.13'
$ws.Range("E34").Value = '  -1.29%  '

$ws.Range("E35").Value = '  -1.05%  '

$ws.Range("E36").Value = '  -0.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.565'
$ws.Range("E37").Value = '  -0.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.880'
$ws.Range("E38").Value = '  -1.12%  '

$ws.Range("E39").Value = '  +0.15%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.897'
$ws.Range("E40").Value = '  +9.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '70.16'
$ws.Range("E41").Value = '  +7.72%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("E43").Value = '  -1.05%  '

$ws.Range("E44").Value = '  +1.20%  '

$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("E46").Value = '  -0.55%  '

$ws.Range("D47").Value = '1.781.20'
$ws.Range("E47").Value = '  -0.49%  '

$ws.Range("E48").Value = '  +3.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.77'
$ws.Range("E49").Value = '  -1.62%  '

$ws.Range("E50").Value = '  -0.70%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0993'
$ws.Range("E51").Value = '  +0.02%  '
